$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item(1)

# Delete row 5 entirely (the "MuSCs" row), shifting nothing up from below
# since it's the last row.
$ws.Rows.Item(5).Delete()

# Row 2 (ECs / Ccl4 / Ackr2 / FAPs) - update numeric values
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.2400906666666667
$ws.Range("H2").Value = 0.720272
$ws.Range("I2").Value = 0.001631540293869566
$ws.Range("J2").Value = 0.001631540293869566
$ws.Range("Q2").Value = 0.01498045714666667
$ws.Range("R2").Value = 0.13482411432
$ws.Range("S2").Value = 0.001631540293869566
$ws.Range("T2").Value = 0.001631540293869566

# Row 3 (Inflammatory-Mac / Ccl4 / Ackr2 / FAPs) - update numeric values
$ws.Range("G3").Value = 91.82408133333333
$ws.Range("H3").Value = 275.472244
$ws.Range("I3").Value = 0.623992138981758
$ws.Range("J3").Value = 0.623992138981758
$ws.Range("Q3").Value = 5.729363554793332
$ws.Range("R3").Value = 51.56427199314
$ws.Range("S3").Value = 0.623992138981758
$ws.Range("T3").Value = 0.623992138981758

# Row 4 was "MuSCs" - now becomes the "Resolving-Mac" row (what used to be
# row 5) with recomputed TPM-based values.
$ws.Range("A4").Value = "Resolving-Mac"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 55.09165833333333
$ws.Range("H4").Value = 165.274975
$ws.Range("I4").Value = 0.3743763207243725
$ws.Range("J4").Value = 0.3743763207243724
$ws.Range("Q4").Value = 3.437444021708334
$ws.Range("R4").Value = 30.936996195375
$ws.Range("S4").Value = 0.3743763207243725
$ws.Range("T4").Value = 0.3743763207243724
